# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# price table with refreshed values/percentages.
#
# Note: Price values are stored as literal text (e.g. "35.608.37",
# "92.10") rather than numbers, matching the source data feed's
# formatting (including any trailing zeros and multi-dot separators).
# A leading apostrophe forces Excel to keep the entry as text instead of
# auto-coercing it to a number (which would silently drop meaningful
# trailing zeros, e.g. "92.10" -> 92.1); the style is then reset back to
# Normal so no stray "quote prefix" cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''35.608.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Value = '''1.983.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.63%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''245.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("E6").Value = '  -4.58%  '
$ws.Range("D7").Value = '''57.78'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.90%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''58.58'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("D12").Value = '''0.104'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.76%  '
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '''14.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '''2.271.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.73%  '
$ws.Range("E16").Value = '  -2.62%  '
$ws.Range("D17").Value = '''1.971.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.28%  '
$ws.Range("D18").Value = '''17.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.64%  '
$ws.Range("D19").Value = '''35.571.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("D20").Value = '''71.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = '''0.0₃0844'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("D22").Value = '''5.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.64%  '
$ws.Range("D23").Value = '''232.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '''2.56'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +19.86%  '
$ws.Range("E26").Value = '  -2.96%  '
$ws.Range("D28").Value = '''9.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").Value = '''19.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.46%  '
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("E31").Value = '  -3.89%  '
$ws.Range("E32").Value = '  -7.49%  '
$ws.Range("D33").Value = '''0.0957'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +16.69%  '
$ws.Range("D34").Value = '''0.0593'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '''2.43'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.15%  '
$ws.Range("D36").Value = '''4.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.77%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '''1.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.52%  '
$ws.Range("D39").Value = '''5.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.30%  '
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("E43").Value = '  -1.91%  '
$ws.Range("D44").Value = '''7.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '''92.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").Value = '''16.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("E47").Value = '  -2.05%  '
$ws.Range("D48").Value = '''1.372.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.48%  '
$ws.Range("D49").Value = '''2.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("D50").Value = '''46.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.86%  '
$ws.Range("E51").Value = '  -1.74%  '
